$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("termsWithMulitpleLabels")
$ws1.Range("D200").Value = "test"
$ws3 = $wb.Worksheets.Item("termWithDifferentParent")
$ws3.Range("D220").Value = "gates_perch"
